# "add food type in this"
# Rename the stray "sss" header to "food Type (halal/haram)" and move it
# from its old stray location (H11) into the new last column (L1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newHeader = "food Type (halal/haram)"

# H11 currently holds the shared string "sss". Updating it in place renames
# that shared string (it is the only cell using it), then we copy the same
# text into the new header cell L1.
$ws.Range("H11").Value = $newHeader
$ws.Range("L1").Value = $newHeader

# Give the new column a reasonable width (closest achievable to 23.41 chars).
$ws.Columns.Item(12).ColumnWidth = 22.5

# Remove the now stray row 11 (shifts the used range back up to row 6).
$ws.Rows.Item(11).Delete()

# Match the author's final selection/active cell.
$ws.Range("J4").Select()
